$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 37038076
$ws.Range("I33").Value = 739.13043
$ws.Range("J33").Value = 250002750
$ws.Range("K33").Value = 739.13043
$ws.Range("L33").Value = 250002750
$ws.Range("M33").Value = -510.13043
$ws.Range("N33").Value = -250003208
$ws.Range("H40").Value = 2506.5557
$ws.Range("I40").Value = 1525.9412
$ws.Range("J40").Value = 3383.9473
$ws.Range("K40").Value = 1525.9412
$ws.Range("L40").Value = 3383.9473
$ws.Range("M40").Value = -1350.9412
$ws.Range("N40").Value = -3733.9473
$ws.Range("H137").Value = 976454.4399999999
$ws.Range("J137").Value = 1236285
$ws.Range("L137").Value = 3708855
$ws.Range("N137").Value = -3713955

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 225.5
$ws.Range("I4").Value = 225.5
$ws.Range("K4").Value = 225.5
$ws.Range("M4").Value = -109.5
$ws.Range("H9").Value = 20861.572
$ws.Range("I9").Value = 70008
$ws.Range("J9").Value = 12670.5
$ws.Range("K9").Value = 70008
$ws.Range("L9").Value = 12670.5
$ws.Range("M9").Value = -69838
$ws.Range("N9").Value = -13010.5
$ws.Range("H20").Value = 20861.572
$ws.Range("I20").Value = 70008
$ws.Range("J20").Value = 12670.5
$ws.Range("K20").Value = 70008
$ws.Range("L20").Value = 12670.5
$ws.Range("M20").Value = -69738
$ws.Range("N20").Value = -13210.5
$ws.Range("H23").Value = 15571.429
$ws.Range("J23").Value = 15571.429
$ws.Range("L23").Value = 15571.429
$ws.Range("N23").Value = -16089.429

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 210
$ws.Range("I11").Value = 210
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 210
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -70
$ws.Range("N11").ClearContents()
$ws.Range("H15").Value = 7800
$ws.Range("J15").Value = 7800
$ws.Range("L15").Value = 7800
$ws.Range("N15").Value = -8254
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H22").Value = 234.11111
$ws.Range("I22").Value = 238.375
$ws.Range("K22").Value = 238.375
$ws.Range("M22").Value = -65.375

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1175.2153
$ws.Range("I31").Value = 1013.6857
$ws.Range("J31").Value = 1363.6666
$ws.Range("K31").Value = 1013.6857
$ws.Range("L31").Value = 1363.6666
$ws.Range("M31").Value = -718.6857
$ws.Range("N31").Value = -1953.6666
$ws.Range("H34").Value = 1175.2153
$ws.Range("I34").Value = 1013.6857
$ws.Range("J34").Value = 1363.6666
$ws.Range("K34").Value = 1013.6857
$ws.Range("L34").Value = 1363.6666
$ws.Range("M34").Value = -811.6857
$ws.Range("N34").Value = -1767.6666
$ws.Range("H41").Value = 9332.666999999999
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H51").Value = 24000
$ws.Range("J51").Value = 24000
$ws.Range("L51").Value = 24000
$ws.Range("N51").Value = -25472
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H61").Value = 24000
$ws.Range("J61").Value = 24000
$ws.Range("L61").Value = 24000
$ws.Range("N61").Value = -24696
$ws.Range("H63").Value = 17925
$ws.Range("I63").Value = 10000
$ws.Range("J63").Value = 25850
$ws.Range("K63").Value = 10000
$ws.Range("L63").Value = 25850
$ws.Range("M63").Value = -9314
$ws.Range("N63").Value = -27222
$ws.Range("H66").Value = 17925
$ws.Range("I66").Value = 10000
$ws.Range("J66").Value = 25850
$ws.Range("K66").Value = 30000
$ws.Range("L66").Value = 77550
$ws.Range("M66").Value = -26568
$ws.Range("N66").Value = -84414
$ws.Range("H81").Value = 29664
$ws.Range("J81").Value = 29664
$ws.Range("L81").Value = 29664
$ws.Range("N81").Value = -31660
$ws.Range("H84").Value = 29664
$ws.Range("J84").Value = 29664
$ws.Range("L84").Value = 88992
$ws.Range("N84").Value = -98976

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1273.5
$ws.Range("I68").Value = 655.2143
$ws.Range("J68").Value = 2355.5
$ws.Range("K68").Value = 1965.6429
$ws.Range("L68").Value = 7066.5
$ws.Range("M68").Value = -1154.6429
$ws.Range("N68").Value = -8688.5
$ws.Range("H71").Value = 1273.5
$ws.Range("I71").Value = 655.2143
$ws.Range("J71").Value = 2355.5
$ws.Range("K71").Value = 5896.928699999999
$ws.Range("L71").Value = 21199.5
$ws.Range("M71").Value = -1840.928699999999
$ws.Range("N71").Value = -29311.5
$ws.Range("H129").Value = 1949.8148
$ws.Range("J129").Value = 2140.611
$ws.Range("L129").Value = 6421.833
$ws.Range("N129").Value = -16421.833
$ws.Range("H131").Value = 908.63
$ws.Range("J131").Value = 918.18555
$ws.Range("L131").Value = 2754.55665
$ws.Range("N131").Value = -12834.55665

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 199.6
$ws.Range("I2").Value = 126.2
$ws.Range("J2").Value = 273
$ws.Range("K2").Value = 126.2
$ws.Range("L2").Value = 273
$ws.Range("M2").Value = -13.2
$ws.Range("N2").Value = -499
$ws.Range("H92").Value = 19987
$ws.Range("J92").Value = 19987
$ws.Range("L92").Value = 19987
$ws.Range("N92").Value = -23731

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 875.1539
$ws.Range("I46").Value = 975
$ws.Range("J46").Value = 857
$ws.Range("K46").Value = 975
$ws.Range("L46").Value = 857
$ws.Range("M46").Value = -787
$ws.Range("N46").Value = -1233
$ws.Range("H55").Value = 606.65
$ws.Range("I55").Value = 195.08333
$ws.Range("J55").Value = 1224
$ws.Range("K55").Value = 195.08333
$ws.Range("L55").Value = 1224
$ws.Range("M55").Value = -22.08332999999999
$ws.Range("N55").Value = -1570
$ws.Range("H62").Value = 29000
$ws.Range("J62").Value = 29000
$ws.Range("L62").Value = 29000
$ws.Range("N62").Value = -30248
$ws.Range("H65").Value = 29000
$ws.Range("J65").Value = 29000
$ws.Range("L65").Value = 87000
$ws.Range("N65").Value = -93240
$ws.Range("H103").Value = 28602
$ws.Range("J103").Value = 28602
$ws.Range("L103").Value = 28602
$ws.Range("N103").Value = -30946

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 10000
$ws.Range("J54").Value = 10000
$ws.Range("L54").Value = 10000
$ws.Range("N54").Value = -11040
$ws.Range("H136").Value = 1447.0513
$ws.Range("I136").Value = 967.75
$ws.Range("J136").Value = 1951.579
$ws.Range("K136").Value = 2903.25
$ws.Range("L136").Value = 5854.737
$ws.Range("M136").Value = -353.25
$ws.Range("N136").Value = -10954.737
